# This script reorders the field-observation rows in the "Artfynd" sheet so
# that each row's full set of values (columns A..AY) matches the row id
# (column A) ordering required by the target revision.
#
# Concretely, the following row-content moves happen (1-based worksheet rows):
#   new row 2  <- old row 3
#   new row 3  <- old row 2
#   new row 4  <- old row 5
#   new row 5  <- old row 6
#   new row 6  <- old row 4
#   new row 8  <- old row 10
#   new row 9  <- old row 8
#   new row 10 <- old row 9
#
# Rows 1 (header), 7 and 11 (and anything beyond AY) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 51  # AY

# A few columns hold values that look numeric/date-like but are authored as
# plain text in this workbook (Antal, Startdatum, Slutdatum). Excel's COM
# layer auto-detects such literals and would silently coerce them into real
# numbers/dates if we just assigned them normally, so for those columns we
# write using a leading apostrophe, which forces a text interpretation -
# exactly like typing '20 or '2023-08-12 directly into Excel.
$forceTextCols = @(9, 25, 27)   # I = Antal, Y = Startdatum, AA = Slutdatum

# Source -> destination row mapping (destination = key, source = value).
$rowSourceForDest = @{
    2  = 3
    3  = 2
    4  = 5
    5  = 6
    6  = 4
    8  = 10
    9  = 8
    10 = 9
}

# 1) Snapshot every value needed as a source, BEFORE any writes happen,
#    so that overlapping/cyclic swaps don't clobber data we still need.
$snapshot = @{}
foreach ($srcRow in ($rowSourceForDest.Values | Sort-Object -Unique)) {
    $rowVals = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value2
    }
    $snapshot[$srcRow] = $rowVals
}

# 2) Write the snapshotted rows into their destination rows.
foreach ($destRow in ($rowSourceForDest.Keys | Sort-Object)) {
    $srcRow = $rowSourceForDest[$destRow]
    $rowVals = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $val = $rowVals[$col]
        $cell = $ws.Cells.Item($destRow, $col)
        if (($forceTextCols -contains $col) -and ($null -ne $val) -and ($val -ne "")) {
            $cell.Value2 = "'" + $val
        } else {
            $cell.Value2 = $val
        }
    }
}
